# Rename the "*img" sheets to "img*"
$wb = $excel.ActiveWorkbook

$renames = @{
    "himg" = "imgh"
    "timg" = "imgt"
    "simg" = "imgs"
    "gimg" = "imgg"
    "wimg" = "imgw"
    "bimg" = "imgb"
    "eimg" = "imge"
}

foreach ($oldName in $renames.Keys) {
    $ws = $wb.Worksheets.Item($oldName)
    $ws.Name = $renames[$oldName]
}

# Activate the last sheet ("imge", previously "eimg") so it becomes the
# selected / active tab, and give it an explicit cell selection.
$lastSheet = $wb.Worksheets.Item("imge")
$lastSheet.Activate()
$lastSheet.Range("A1").Select()

# Touch the alignment/protection of "misc"!E2 so Excel writes the
# applyAlignment / applyProtection flags on its style (re-applying the
# same visual values: general horizontal alignment, locked cell).
$miscSheet = $wb.Worksheets.Item("misc")
$cell = $miscSheet.Range("E2")
$cell.HorizontalAlignment = -4128
$cell.Locked = $true
